# Update gh-pages to output generated at 456a3b4
# Refreshes the "想去人数" (interest count, column F) and "最低票价"
# (minimum price, column G) figures on the 展览 and 全部类型 sheets.

$wb = $excel.ActiveWorkbook

# Row -> new F-column value, for rows that sit at the SAME row number on
# both the 展览 and 全部类型 sheets (everything above the extra row that
# 全部类型 has inserted at its row 31).
$commonUpdates = @{
    5  = 122
    7  = 92
    8  = 1213
    9  = 17189
    10 = 312
    11 = 219
    13 = 6570
    14 = 670
    15 = 142
    17 = 44
    18 = 134
    20 = 133
    21 = 47
    24 = 19
    25 = 27
    26 = 251
    28 = 87
    29 = 5101
    30 = 522
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Row 3: F3 ticks up to 1162, and G3 flips from a numeric price (50) to
    # the inline string "不可售" (not for sale).
    $ws.Range("F3").Value = 1162
    $ws.Range("G3").Value = "不可售"

    foreach ($row in $commonUpdates.Keys) {
        $ws.Cells.Item($row, 6).Value = $commonUpdates[$row]
    }
}

# 展览 rows 31-39 (no offset).
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Cells.Item(31, 6).Value = 49
$wsExpo.Cells.Item(32, 6).Value = 11665
$wsExpo.Cells.Item(35, 6).Value = 177
$wsExpo.Cells.Item(36, 6).Value = 245
$wsExpo.Cells.Item(37, 6).Value = 3877
$wsExpo.Cells.Item(39, 6).Value = 81

# 全部类型 carries one extra row above this point, so its matching rows are
# shifted down by one relative to 展览.
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Cells.Item(32, 6).Value = 49
$wsAll.Cells.Item(33, 6).Value = 11665
$wsAll.Cells.Item(36, 6).Value = 177
$wsAll.Cells.Item(37, 6).Value = 245
$wsAll.Cells.Item(38, 6).Value = 3877
$wsAll.Cells.Item(40, 6).Value = 81

Write-Output "Applied F/G column refresh to 展览 and 全部类型"
